$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated average runtimes for days 1-13, plus new rows for days 14-16 (solved day 16 part 2)
$values = @(
    0.01610408,
    0.031768100000000001,
    0.03043978,
    0.0065995200000000002,
    0.010931099999999999,
    0.04096764,
    0.02510078,
    0.0221143,
    0.035815520000000003,
    0.1863002,
    0.14555241999999999,
    0.136183,
    0.0042554000000000003,
    0.14823896,
    0,
    1.67006312
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $day = $i + 1
    $row = $i + 3
    $ws.Cells.Item($row, 1).Value = $day
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$ws.Range("B22").Select()
